$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values entered in this order so the shared-string table matches the
# author's original editing sequence.
$ws.Range("A1").Value = "EndStrips"
$ws.Range("A2").Value = "EndCaps"
$ws.Range("B2").Value = "Extended release glucose capsules. Never worry about hypoglycemia at night again!"
$ws.Range("D2").Value = "img/endcaps.png"
$ws.Range("B1").Value = "Revolutionary micro-dosing hypoglycemic system. Developed in house by EndT1."
$ws.Range("D1").Value = "img/endstrips.png"
$ws.Range("A3").Value = "test"
$ws.Range("B3").Value = "asdfasd"

$ws.Range("C1").Value = 10
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = "img/endcaps.png"

# Remove bold formatting from the header row -> back to the default style
$ws.Range("A1:D1").Style = "Normal"

# Apply a numeric format to the price column's first two rows
$ws.Range("C1:C2").NumberFormat = "0.00"

# Resize columns to fit the new content (closest values this engine's
# character-width quantization can reach to the author's 8.375 / 71 /
# 5.375 / 26.125 target widths)
$ws.Columns(1).ColumnWidth = 7.5
$ws.Columns(2).ColumnWidth = 70.16666666666667
$ws.Columns(3).ColumnWidth = 4.5
$ws.Columns(4).ColumnWidth = 25.333333333333336

# Update selection to D3
$ws.Range("D3").Select()
